$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.737.80'
$ws.Range('E2').Value = '  +2.43%  '

$ws.Range('D3').Value = '3.049.30'
$ws.Range('E3').Value = '  +2.03%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.52%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.17%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  +5.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.62'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.09%  '

$ws.Range('E10').Value = '  +7.43%  '

$ws.Range('E11').Value = '  +5.45%  '

$ws.Range('E12').Value = '  +2.19%  '

$ws.Range('D13').Value = '3.576.85'
$ws.Range('E13').Value = '  +2.18%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.97'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.78%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000170'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +16.65%  '

$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.60%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '57.725.62'
$ws.Range('E17').Value = '  +2.44%  '

$ws.Range('D18').Value = '3.054.99'
$ws.Range('E18').Value = '  +2.20%  '

$ws.Range('E19').Value = '  +5.70%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '338.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.06%  '

$ws.Range('E22').Value = '  +0.03%  '

$ws.Range('E23').Value = '  +7.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.94%  '

$ws.Range('E25').Value = '  +6.82%  '

$ws.Range('D26').Value = '0.0₃0980'
$ws.Range('E26').Value = '  +9.15%  '

$ws.Range('E27').Value = '  +0.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.86%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.43%  '

$ws.Range('E31').Value = '  +5.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.38%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '156.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.08%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.14%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.99'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.57%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.27%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.09'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.59%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0706'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.65%  '

$ws.Range('D39').Value = '3.086.24'
$ws.Range('E39').Value = '  +2.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.27%  '

$ws.Range('E41').Value = '  +9.02%  '

$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('E43').Value = '  +5.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.663'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.69%  '

$ws.Range('D45').Value = '2.329.17'
$ws.Range('E45').Value = '  +4.24%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.13%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.05%  '

$ws.Range('E48').Value = '  +4.33%  '

$ws.Range('E49').Value = '  +4.24%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.79%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0899'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.12%  '
